# Updated cryptos list on Mon Apr 29 08:59:18 UTC 2024 with GitHub Actions
# Applies per-cell price/volume(1h) updates (and a few row label/link swaps)
# to the crypto ranking sheet, matching the upstream data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell values in column D that look like plain numbers must be forced back to
# text (a leading apostrophe / quote-prefix) so Excel does not silently convert
# them to numeric values and drop formatting such as trailing zeros.

$ws.Range("D2").Value = "62.491.36"
$ws.Range("E2").Value = "  -2.07%  "
$ws.Range("D3").Value = "3.171.84"
$ws.Range("E3").Value = "  -4.24%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'586.47"
$ws.Range("E5").Value = "  -2.51%  "
$ws.Range("D6").Value = "'134.80"
$ws.Range("E6").Value = "  -5.57%  "
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("D8").Value = "3.170.58"
$ws.Range("E8").Value = "  -4.24%  "
$ws.Range("D9").Value = "'0.503"
$ws.Range("E9").Value = "  -3.65%  "
$ws.Range("E10").Value = "  -5.47%  "
$ws.Range("D11").Value = "'5.24"
$ws.Range("E11").Value = "  -5.53%  "
$ws.Range("D12").Value = "'0.451"
$ws.Range("E12").Value = "  -4.85%  "
$ws.Range("D13").Value = "'0.0000234"
$ws.Range("E13").Value = "  -5.96%  "
$ws.Range("D14").Value = "'33.06"
$ws.Range("E14").Value = "  -4.99%  "
$ws.Range("D15").Value = "3.694.75"
$ws.Range("E15").Value = "  -4.24%  "
$ws.Range("E16").Value = "  -1.75%  "
$ws.Range("D17").Value = "3.167.17"
$ws.Range("E17").Value = "  -4.34%  "
$ws.Range("D18").Value = "62.459.12"
$ws.Range("E18").Value = "  -2.25%  "
$ws.Range("E19").Value = "  -5.39%  "
$ws.Range("D20").Value = "'455.48"
$ws.Range("E20").Value = "  -5.43%  "
$ws.Range("D21").Value = "'13.84"
$ws.Range("E21").Value = "  -2.69%  "
$ws.Range("D22").Value = "'0.702"
$ws.Range("E22").Value = "  -4.59%  "
$ws.Range("D23").Value = "'7.61"
$ws.Range("E23").Value = "  -4.96%  "
$ws.Range("D24").Value = "'13.33"
$ws.Range("E24").Value = "  -1.64%  "
$ws.Range("D25").Value = "'82.43"
$ws.Range("E25").Value = "  -2.77%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  -0.12%  "
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("E28").Value = "  -3.44%  "
$ws.Range("D29").Value = "'6.89"
$ws.Range("E29").Value = "  -5.38%  "
$ws.Range("D30").Value = "'7.81"
$ws.Range("E30").Value = "  -4.33%  "
$ws.Range("E31").Value = "  -7.28%  "
$ws.Range("D32").Value = "'27.23"
$ws.Range("E32").Value = "  -7.14%  "
$ws.Range("E33").Value = "  -3.36%  "
$ws.Range("D34").Value = "'2.39"
$ws.Range("E34").Value = "  -6.32%  "
$ws.Range("E35").Value = "  -6.71%  "
$ws.Range("D36").Value = "'5.78"
$ws.Range("E36").Value = "  -3.27%  "
$ws.Range("D37").Value = "'51.12"
$ws.Range("E37").Value = "  -3.56%  "
$ws.Range("D38").Value = "0.0₃0693"
$ws.Range("E38").Value = "  -8.03%  "
$ws.Range("D39").Value = "'0.0384"
$ws.Range("E39").Value = "  -4.77%  "
$ws.Range("D40").Value = "'411.72"
$ws.Range("E40").Value = "  -4.35%  "
$ws.Range("D41").Value = "2.909.28"
$ws.Range("E41").Value = "  -4.59%  "
$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D42").Value = "'2.66"
$ws.Range("E42").Value = "  -3.35%  "
$ws.Range("B43").Value = "Cosmos"
$ws.Range("C43").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D43").Value = "'7.98"
$ws.Range("E43").Value = "  -5.25%  "
$ws.Range("B44").Value = "Kaspa"
$ws.Range("C44").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D44").Value = "'0.112"
$ws.Range("E44").Value = "  -0.30%  "
$ws.Range("D45").Value = "'0.249"
$ws.Range("E45").Value = "  -6.62%  "
$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D46").Value = "'2.13"
$ws.Range("E46").Value = "  -3.37%  "
$ws.Range("B47").Value = "USDe"
$ws.Range("C47").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D47").Value = "'0.999"
$ws.Range("E47").Value = "  -0.05%  "
$ws.Range("D48").Value = "'35.83"
$ws.Range("E48").Value = "  -0.92%  "
$ws.Range("D49").Value = "'124.72"
$ws.Range("E49").Value = "  +0.36%  "
$ws.Range("D50").Value = "'25.30"
$ws.Range("E50").Value = "  -4.46%  "
$ws.Range("E51").Value = "  -3.92%  "
